$wb = $excel.ActiveWorkbook


# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 4465291.5
$ws.Range("I129").Value = 50001456
$ws.Range("J129").Value = 961.7646999999999
$ws.Range("K129").Value = 150004368
$ws.Range("L129").Value = 2885.2941
$ws.Range("M129").Value = -149999368
$ws.Range("N129").Value = -12885.2941
$ws.Range("H137").Value = 2265.4465
$ws.Range("I137").Value = 2259.2046
$ws.Range("J137").Value = 2288.3333
$ws.Range("K137").Value = 6777.6138
$ws.Range("L137").Value = 6864.999899999999
$ws.Range("M137").Value = -4227.6138
$ws.Range("N137").Value = -11964.9999
$ws.Range("H138").Value = 4709.953
$ws.Range("J138").Value = 6016.911
$ws.Range("L138").Value = 18050.733
$ws.Range("N138").Value = -28330.733

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1004.8
$ws.Range("I74").Value = 909.6
$ws.Range("J74").Value = 1100
$ws.Range("K74").Value = 909.6
$ws.Range("L74").Value = 1100
$ws.Range("M74").Value = -35.60000000000002
$ws.Range("N74").Value = -2848
$ws.Range("H77").Value = 1004.8
$ws.Range("I77").Value = 909.6
$ws.Range("J77").Value = 1100
$ws.Range("K77").Value = 4548
$ws.Range("L77").Value = 5500
$ws.Range("M77").Value = -180
$ws.Range("N77").Value = -14236
$ws.Range("H88").Value = 1579.8
$ws.Range("I88").Value = 1624.75
$ws.Range("K88").Value = 1624.75
$ws.Range("M88").Value = -1218.75
$ws.Range("H91").Value = 1579.8
$ws.Range("I91").Value = 1624.75
$ws.Range("K91").Value = 1624.75
$ws.Range("M91").Value = -220.75
$ws.Range("H134").Value = 29950
$ws.Range("J134").Value = 29950
$ws.Range("L134").Value = 29950
$ws.Range("N134").Value = -40090
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 21450.154
$ws.Range("I86").Value = 1474.25
$ws.Range("K86").Value = 1474.25
$ws.Range("M86").Value = -351.25
$ws.Range("H89").Value = 21450.154
$ws.Range("I89").Value = 1474.25
$ws.Range("K89").Value = 7371.25
$ws.Range("M89").Value = -1755.25
$ws.Range("H141").Value = 37045.453
$ws.Range("J141").Value = 29642.857
$ws.Range("L141").Value = 29642.857
$ws.Range("N141").Value = -40002.857

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 53126.75
$ws.Range("I4").Value = 2501
$ws.Range("J4").Value = 70002
$ws.Range("K4").Value = 2501
$ws.Range("L4").Value = 70002
$ws.Range("M4").Value = -2389
$ws.Range("N4").Value = -70226
$ws.Range("H31").Value = 2657.4119
$ws.Range("I31").Value = 1723.2609
$ws.Range("J31").Value = 3424.75
$ws.Range("K31").Value = 1723.2609
$ws.Range("L31").Value = 3424.75
$ws.Range("M31").Value = -1428.2609
$ws.Range("N31").Value = -4014.75
$ws.Range("H34").Value = 2657.4119
$ws.Range("I34").Value = 1723.2609
$ws.Range("J34").Value = 3424.75
$ws.Range("K34").Value = 1723.2609
$ws.Range("L34").Value = 3424.75
$ws.Range("M34").Value = -1521.2609
$ws.Range("N34").Value = -3828.75
$ws.Range("H52").Value = 30000
$ws.Range("J52").Value = 30000
$ws.Range("L52").Value = 30000
$ws.Range("N52").Value = -30588
$ws.Range("H58").Value = 8476609
$ws.Range("I58").Value = 1064.0465
$ws.Range("J58").Value = 31254634
$ws.Range("K58").Value = 1064.0465
$ws.Range("L58").Value = 31254634
$ws.Range("M58").Value = -861.0464999999999
$ws.Range("N58").Value = -31255040
$ws.Range("H62").Value = 4509.4
$ws.Range("J62").Value = 4585.5
$ws.Range("L62").Value = 4585.5
$ws.Range("N62").Value = -5833.5
$ws.Range("H65").Value = 4509.4
$ws.Range("J65").Value = 4585.5
$ws.Range("L65").Value = 22927.5
$ws.Range("N65").Value = -29167.5
$ws.Range("H127").Value = 33000
$ws.Range("J127").Value = 33000
$ws.Range("L127").Value = 33000
$ws.Range("N127").Value = -42920
$ws.Range("H136").Value = 8476609
$ws.Range("I136").Value = 1064.0465
$ws.Range("J136").Value = 31254634
$ws.Range("K136").Value = 3192.1395
$ws.Range("L136").Value = 93763902
$ws.Range("M136").Value = -642.1394999999998
$ws.Range("N136").Value = -93769002

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 894.5714
$ws.Range("I4").Value = 86
$ws.Range("J4").Value = 1501
$ws.Range("K4").Value = 258
$ws.Range("L4").Value = 4503
$ws.Range("M4").Value = -146
$ws.Range("N4").Value = -4727
$ws.Range("H58").Value = 955.55554
$ws.Range("I58").Value = 483.33334
$ws.Range("J58").Value = 1900
$ws.Range("K58").Value = 1450.00002
$ws.Range("L58").Value = 5700
$ws.Range("M58").Value = -1322.00002
$ws.Range("N58").Value = -5956
$ws.Range("H92").Value = 1972
$ws.Range("I92").Value = 486.33334
$ws.Range("J92").Value = 2343.4167
$ws.Range("K92").Value = 1459.00002
$ws.Range("L92").Value = 7030.250100000001
$ws.Range("M92").Value = -211.0000199999999
$ws.Range("N92").Value = -9526.250100000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 31.619047
$ws.Range("I2").Value = 29.071428
$ws.Range("J2").Value = 36.714287
$ws.Range("K2").Value = 29.071428
$ws.Range("L2").Value = 36.714287
$ws.Range("M2").Value = 83.928572
$ws.Range("N2").Value = -262.714287
$ws.Range("H5").Value = 13666.333
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 13666.333
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 13666.333
$ws.Range("M5").ClearContents()
$ws.Range("H137").Value = 29583.334
$ws.Range("J137").Value = 29583.334
$ws.Range("L137").Value = 29583.334
$ws.Range("N137").Value = -39783.334

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 242960.12
$ws.Range("J2").Value = 43769.46
$ws.Range("L2").Value = 43769.46
$ws.Range("N2").Value = -43993.46
$ws.Range("H22").Value = 200003260
$ws.Range("I22").Value = 500000100
$ws.Range("J22").Value = 5383.6665
$ws.Range("K22").Value = 500000100
$ws.Range("L22").Value = 5383.6665
$ws.Range("M22").Value = -499999805
$ws.Range("N22").Value = -5973.6665
$ws.Range("H27").Value = 200003260
$ws.Range("I27").Value = 500000100
$ws.Range("J27").Value = 5383.6665
$ws.Range("K27").Value = 500000100
$ws.Range("L27").Value = 5383.6665
$ws.Range("M27").Value = -499999993
$ws.Range("N27").Value = -5597.6665
$ws.Range("H46").Value = 1602.875
$ws.Range("I46").Value = 402.2857
$ws.Range("J46").Value = 2536.6667
$ws.Range("K46").Value = 402.2857
$ws.Range("L46").Value = 2536.6667
$ws.Range("M46").Value = -214.2857
$ws.Range("N46").Value = -2912.6667
$ws.Range("H55").Value = 841.7368
$ws.Range("I55").Value = 217.28572
$ws.Range("J55").Value = 1206
$ws.Range("K55").Value = 217.28572
$ws.Range("L55").Value = 1206
$ws.Range("M55").Value = -44.28572
$ws.Range("N55").Value = -1552
$ws.Range("H68").Value = 2886.8
$ws.Range("I68").Value = 1166.5555
$ws.Range("J68").Value = 5467.1665
$ws.Range("K68").Value = 1166.5555
$ws.Range("L68").Value = 5467.1665
$ws.Range("M68").Value = -417.5554999999999
$ws.Range("N68").Value = -6965.1665
$ws.Range("H71").Value = 2886.8
$ws.Range("I71").Value = 1166.5555
$ws.Range("J71").Value = 5467.1665
$ws.Range("K71").Value = 5832.7775
$ws.Range("L71").Value = 27335.8325
$ws.Range("M71").Value = -2088.7775
$ws.Range("N71").Value = -34823.8325
$ws.Range("H135").Value = 29775.572
$ws.Range("J135").Value = 29775.572
$ws.Range("L135").Value = 29775.572
$ws.Range("N135").Value = -39915.572
